$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(113, 105, 8, 35024.97, 348, 1397.901641398479, 92.92, 7.83, 3125.01, 321.2, -55.89, 188.75, 1.67, 33279.07, 20.06)
    3 = @(111, 106, 5, 34537.06999999999, 221, 1353.799795035902, 95.5, 7.63, 3125.01, 313.61, -56.04, 187.25, 1.69, 32962.27, 21.93)
    4 = @(107, 105, 2, 34071.07000000001, 91, 1311.542161780268, 98.13, 7.61, 3125.01, 312.27, -59.91, 182.25, 1.7, 32668.53, 24.29)
    5 = @(109, 106, 3, 34537.06999999999, 121, 1338.40762331873, 97.25, 7.63, 3125.01, 313.61, -54.93, 185.75, 1.7, 33077.66, 23.67)
    6 = @(107, 104, 3, 33936.07, 121, 1320.22006467283, 97.2, 7.66, 3125.01, 314.04, -54.93, 183, 1.71, 32494.85, 23.55)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$wb.Save()
